# regen sval data to filter save games
# Overwrite the computed stat columns (B:G) for rows 2-15 on the active sheet
# with the regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
    3  = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 6.348428708163715)
    4  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    5  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    6  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    7  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    8  = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 10.19245300693656, 1, 16.84135478251809)
    9  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
    10 = @(1.455362044514542, 117.745847958593, 0.1494219747398047, 10.19245300693656, 1, 129.5430849847839)
    11 = @(0.1190320826869504, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 1.069511820747223)
    12 = @(1.455362044514542, 10.34677158129881, 22.3905356188092, 10.19245300693656, 1, 44.38512225155911)
    13 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 3.754798637575387)
    14 = @(0.6606524410359556, 0.306821227259698, 22.3905356188092, 0.4942365360607697, 1, 23.85224582316562)
    15 = @(0.2917716402565462, 1.655778082260271, 3.537761648806719, 10.19245300693656, 1, 15.67776437826009)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
